$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

    $ws.Range("D2").Value = 0.85
    $ws.Range("E2").Value = 0.6123369849394824
    $ws.Range("G2").Value = 0.7874469375379018
    $ws.Range("H2").Value = 0.7394957983193278
    $ws.Range("E5").Value = 0.7223209035150899
    $ws.Range("F5").Value = 0.8606060606060606
    $ws.Range("G5").Value = 0.8660551849605822
    $ws.Range("H5").Value = 0.8413793103448276
    $ws.Range("E9").Value = 0.7033000723291813
    $ws.Range("F9").Value = 0.8545454545454545
    $ws.Range("G9").Value = 0.8543056397816858
    $ws.Range("H9").Value = 0.8285714285714286
    $ws.Range("D11").Value = 0.73
    $ws.Range("D15").Value = 0.9399999999999998
    $ws.Range("E15").Value = 0.6598332842154305
    $ws.Range("F15").Value = 0.8242424242424242
    $ws.Range("G15").Value = 0.8351273499090358
    $ws.Range("H15").Value = 0.8079470198675497
    $ws.Range("D16").Value = 0.8400000000000001
    $ws.Range("E16").Value = 0.7478564253481768
    $ws.Range("F16").Value = 0.896969696969697
    $ws.Range("G16").Value = 0.8863636363636364
    $ws.Range("H16").Value = 0.8172043010752688
    $ws.Range("D18").Value = 0.97
    $ws.Range("E18").Value = 0.7255281112076279
    $ws.Range("F18").Value = 0.8787878787878788
    $ws.Range("G18").Value = 0.8884297520661157
    $ws.Range("H18").Value = 0.8
    $ws.Range("E19").Value = 0.8548859804241099
    $ws.Range("F19").Value = 0.9393939393939394
    $ws.Range("G19").Value = 0.9442148760330579
    $ws.Range("H19").Value = 0.8936170212765958
    $ws.Range("E21").Value = 0.7865178255868005
    $ws.Range("F21").Value = 0.9151515151515152
    $ws.Range("G21").Value = 0.8987603305785123
    $ws.Range("H21").Value = 0.8444444444444444
    $ws.Range("E22").Value = 0.4806580133202482
    $ws.Range("F22").Value = 0.7454545454545455
    $ws.Range("G22").Value = 0.7685950413223142
    $ws.Range("H22").Value = 0.6315789473684209
    $ws.Range("D23").Value = 0.9299999999999999
    $ws.Range("E23").Value = 0.6640516622105241
    $ws.Range("F23").Value = 0.8424242424242424
    $ws.Range("G23").Value = 0.8636363636363636
    $ws.Range("H23").Value = 0.7547169811320754
    $ws.Range("E24").Value = 0.6645327130973755
    $ws.Range("F24").Value = 0.8303030303030303
    $ws.Range("G24").Value = 0.8698347107438017
    $ws.Range("H24").Value = 0.75
    $ws.Range("D25").Value = 0.991
    $ws.Range("E25").Value = 0.444731204956186
    $ws.Range("F25").Value = 0.7272727272727273
    $ws.Range("G25").Value = 0.7489669421487604
    $ws.Range("H25").Value = 0.6086956521739131
    $ws.Range("E26").Value = 0.5915790933172286
    $ws.Range("F26").Value = 0.7939393939393939
    $ws.Range("G26").Value = 0.8305785123966942
    $ws.Range("H26").Value = 0.7017543859649122
    $ws.Range("E27").Value = 0.6676402299942533
    $ws.Range("F27").Value = 0.8545454545454545
    $ws.Range("G27").Value = 0.8574380165289257
    $ws.Range("H27").Value = 0.76
    $ws.Range("E28").Value = 0.8867495623482351
    $ws.Range("F28").Value = 0.9696969696969697
    $ws.Range("G28").Value = 0.9825174825174825
    $ws.Range("H28").Value = 0.8979591836734693
    $ws.Range("D29").Value = 0.7999999999999999
    $ws.Range("D31").Value = 0.86
    $ws.Range("E31").Value = 0.8676100810544175
    $ws.Range("F31").Value = 0.9636363636363636
    $ws.Range("G31").Value = 0.979020979020979
    $ws.Range("H31").Value = 0.88
    $ws.Range("D33").Value = 0.8099999999999999
    $ws.Range("D35").Value = 0.73
    $ws.Range("E35").Value = 0.8494029733528037
    $ws.Range("F35").Value = 0.9575757575757575
    $ws.Range("G35").Value = 0.9755244755244755
    $ws.Range("H35").Value = 0.8627450980392156
    $ws.Range("D36").Value = 0.76
    $ws.Range("E36").Value = 0.8867495623482351
    $ws.Range("F36").Value = 0.9696969696969697
    $ws.Range("G36").Value = 0.9825174825174825
    $ws.Range("H36").Value = 0.8979591836734693
    $ws.Range("D37").Value = 0.7799999999999999
    $ws.Range("D38").Value = 0.8099999999999999
    $ws.Range("D39").Value = 0.9199999999999999
    $ws.Range("E39").Value = 0.8579757250224321
    $ws.Range("F39").Value = 0.9636363636363636
    $ws.Range("G39").Value = 0.9597902097902098
    $ws.Range("H39").Value = 0.875
    $ws.Range("D40").Value = 0.993
    $ws.Range("E40").Value = 0.8044335210322532
    $ws.Range("F40").Value = 0.9696969696969697
    $ws.Range("G40").Value = 0.8633333333333333
    $ws.Range("H40").Value = 0.8148148148148148
    $ws.Range("D42").Value = 0.999
    $ws.Range("E42").Value = 0.7595545253127499
    $ws.Range("F42").Value = 0.9636363636363636
    $ws.Range("G42").Value = 0.8
    $ws.Range("H42").Value = 0.7499999999999999
    $ws.Range("D52").Value = 0.7600000000000001
    $ws.Range("E59").Value = 0.4672567184860485
    $ws.Range("F59").Value = 0.9454545454545454
    $ws.Range("G59").Value = 0.8113207547169812
    $ws.Range("H59").Value = 0.4705882352941177
    $ws.Range("D60").Value = 0.986
    $ws.Range("D61").Value = 0.9879999999999999
    $ws.Range("E61").Value = 0.700528900717694
    $ws.Range("F61").Value = 0.9818181818181818
    $ws.Range("G61").Value = 0.75
    $ws.Range("H61").Value = 0.6666666666666666
    $ws.Range("D62").Value = 0.9869999999999999
    $ws.Range("E62").Value = 0.700528900717694
    $ws.Range("F62").Value = 0.9818181818181818
    $ws.Range("G62").Value = 0.75
    $ws.Range("H62").Value = 0.6666666666666666
    $ws.Range("D63").Value = 0.97
    $ws.Range("E63").Value = 0.5322721853011699
    $ws.Range("F63").Value = 0.9696969696969697
    $ws.Range("G63").Value = 0.7437106918238994
    $ws.Range("H63").Value = 0.5454545454545454
    $ws.Range("D64").Value = 0.9800000000000001
    $ws.Range("D65").Value = 0.991
    $ws.Range("E65").Value = 0.6009021890119928
    $ws.Range("F65").Value = 0.9757575757575757
    $ws.Range("G65").Value = 0.7468553459119497
    $ws.Range("H65").Value = 0.6
    $ws.Range("D66").Value = 0.991
    $ws.Range("D67").Value = 0.8099999999999999
    $ws.Range("D68").Value = 0.8799999999999999
    $ws.Range("D69").Value = 0.991
    $ws.Range("E69").Value = 0.6009021890119928
    $ws.Range("F69").Value = 0.9757575757575757
    $ws.Range("G69").Value = 0.7468553459119497
    $ws.Range("H69").Value = 0.6
    $ws.Range("D70").Value = 0.9879999999999999
    $ws.Range("D71").Value = 0.983
    $ws.Range("D72").Value = 0.9800000000000001
    $ws.Range("E72").Value = 0.6009021890119928
    $ws.Range("F72").Value = 0.9757575757575757
    $ws.Range("G72").Value = 0.7468553459119497
    $ws.Range("H72").Value = 0.6
    $ws.Range("D80").Value = 0.77
    $ws.Range("E80").Value = 0.6465210112513111
    $ws.Range("F80").Value = 0.9757575757575757
    $ws.Range("G80").Value = 0.9876543209876543
    $ws.Range("H80").Value = 0.6
    $ws.Range("D85").Value = 0.9099999999999998
    $ws.Range("E85").Value = 0.5051814855409226
    $ws.Range("F85").Value = 0.9757575757575757
    $ws.Range("G85").Value = 0.8240740740740741
    $ws.Range("H85").Value = 0.5
    $ws.Range("D86").Value = 0.9199999999999999
    $ws.Range("E86").Value = 0.5051814855409226
    $ws.Range("F86").Value = 0.9757575757575757
    $ws.Range("G86").Value = 0.8240740740740741
    $ws.Range("H86").Value = 0.5
    $ws.Range("D89").Value = 0.9099999999999999
    $ws.Range("E89").Value = 0.5875
    $ws.Range("F89").Value = 0.9757575757575757
    $ws.Range("G89").Value = 0.79375
    $ws.Range("H89").Value = 0.6
    $ws.Range("D93").Value = 0.985
    $ws.Range("D95").Value = 0.7
    $ws.Range("D97").Value = 0.95
    $ws.Range("D100").Value = 0.8899999999999999
    $ws.Range("E100").Value = 0.4732367924632347
    $ws.Range("F100").Value = 0.9393939393939394
    $ws.Range("G100").Value = 0.871875
    $ws.Range("H100").Value = 0.4444444444444444